# Rolling the income-statement table forward by one fiscal year:
# drop the 1396/12 column, shift 1397..1400 left into D..G,
# and populate the new rightmost column H with the 1401/12 figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: financial-period headers ---
$ws.Range("D8").Value = '12 ماهه منتهی به 1397/12'
$ws.Range("E8").Value = '12 ماهه منتهی به 1398/12'
$ws.Range("F8").Value = '12 ماهه منتهی به 1399/12'
$ws.Range("G8").Value = '12 ماهه منتهی به 1400/12'
$ws.Range("H8").Value = '12 ماهه منتهی به 1401/12'

# --- Row 9: publish dates for each period ---
$ws.Range("D9").Value = '1399-03-20 (8)'
$ws.Range("E9").Value = '1400-04-02 (8)'
$ws.Range("F9").Value = '1401-04-08 (8)'
$ws.Range("G9").Value = '1402-02-30 (8)'
$ws.Range("H9").Value = '1402-02-30'

# --- Rows 11-27: shift each metric one period left, append new 1401/12 values in column H ---
# row 11
$ws.Range("D11").Value = 73845433
$ws.Range("E11").Value = 78031948
$ws.Range("F11").Value = 143234768
$ws.Range("G11").Value = 192628444
$ws.Range("H11").Value = 214213606

# row 12
$ws.Range("D12").Value = -29833190
$ws.Range("E12").Value = -35383558
$ws.Range("F12").Value = -61344224
$ws.Range("G12").Value = -146246354
$ws.Range("H12").Value = -147350610

# row 13
$ws.Range("D13").Value = 44012243
$ws.Range("E13").Value = 42648390
$ws.Range("F13").Value = 81890544
$ws.Range("G13").Value = 46382090
$ws.Range("H13").Value = 66862996

# row 14
$ws.Range("D14").Value = -10499423
$ws.Range("E14").Value = -16575782
$ws.Range("F14").Value = -34001119
$ws.Range("G14").Value = -32817902
$ws.Range("H14").Value = -46301021

# row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# row 16
$ws.Range("D16").Value = 3189587
$ws.Range("E16").Value = 3774048
$ws.Range("F16").Value = 9770789
$ws.Range("G16").Value = -689512
$ws.Range("H16").Value = 10839963

# row 17
$ws.Range("D17").Value = 36702407
$ws.Range("E17").Value = 29846656
$ws.Range("F17").Value = 57660214
$ws.Range("G17").Value = 12874676
$ws.Range("H17").Value = 31401938

# row 18
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0

# row 19
$ws.Range("D19").Value = 793837
$ws.Range("E19").Value = -450819
$ws.Range("F19").Value = 247268
$ws.Range("G19").Value = 3737347
$ws.Range("H19").Value = 220497

# row 20
$ws.Range("D20").Value = 37496244
$ws.Range("E20").Value = 29395837
$ws.Range("F20").Value = 57907482
$ws.Range("G20").Value = 16612023
$ws.Range("H20").Value = 31622435

# row 21
$ws.Range("D21").Value = -196171
$ws.Range("E21").Value = -7502
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0

# row 22
$ws.Range("D22").Value = 37300073
$ws.Range("E22").Value = 29388335
$ws.Range("F22").Value = 57907482
$ws.Range("G22").Value = 16612023
$ws.Range("H22").Value = 31622435

# row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# row 24
$ws.Range("D24").Value = 37300073
$ws.Range("E24").Value = 29388335
$ws.Range("F24").Value = 57907482
$ws.Range("G24").Value = 16612023
$ws.Range("H24").Value = 31622435

# row 25
$ws.Range("D25").Value = 15542
$ws.Range("E25").Value = 12245
$ws.Range("F25").Value = 24128
$ws.Range("G25").Value = 6922
$ws.Range("H25").Value = 13176

# row 26
$ws.Range("D26").Value = 2400000
$ws.Range("E26").Value = 2400000
$ws.Range("F26").Value = 2400000
$ws.Range("G26").Value = 2400000
$ws.Range("H26").Value = 2400000

# row 27
$ws.Range("D27").Value = 15542
$ws.Range("E27").Value = 12245
$ws.Range("F27").Value = 24128
$ws.Range("G27").Value = 6922
$ws.Range("H27").Value = 13176

